# feat: Install the dompdf library and create a view for generating reports
# Updates the "Nota" explanatory box on the products import template:
#  - existing note becomes "Nota 1: ..." (unchanged body text)
#  - a new "Nota 2: ..." box is added (merged E6:F6) explaining that example
#    rows should be copied starting at row 18
#  - row heights for rows 6/7 grow to fit the new wrapped text
#  - the sheet selection/scroll position is moved down to the new content

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Nota 1 (was just "Nota: ") in the existing E4:F5 merged box ---
$nota1Bold = "Nota 1: "
$nota1Rest = 'Si en los campos de imágenes usted proporciona urls de las cuales éstas no se pueden descargar, el producto se colocará en estado "inactivo" y en su lugar se agregararán imágenes por defecto.'
$ws.Range("E4").Value = $nota1Bold + $nota1Rest
$ws.Range("E4").Characters(1, $nota1Bold.Length).Font.Bold = $true
$ws.Range("E4").Characters($nota1Bold.Length + 1, $nota1Rest.Length).Font.Bold = $false

# --- Nota 2 (new) in a new E6:F6 merged box ---
$ws.Range("E6:F6").Merge()
$nota2Bold = "Nota 2:"
$nota2Rest = " copie sus productos para importar, siempre desde la fila 18, donde se encuentra actualmente el registro de ejemplo"
$ws.Range("E6").Value = $nota2Bold + $nota2Rest
$ws.Range("E4").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E6").Characters(1, $nota2Bold.Length).Font.Bold = $true
$ws.Range("E6").Characters($nota2Bold.Length + 1, $nota2Rest.Length).Font.Bold = $false

# --- row heights for the taller notes ---
$ws.Rows(6).RowHeight = 39.75
$ws.Rows(7).RowHeight = 27

# --- view / selection moved to where the user was last working ---
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("E11").Select()

Write-Output "done"
